$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.107.27'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.69%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.815.12'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('E4').Value = '  +0.62%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.46'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.613'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.25'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.324'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0684'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.12%  '
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.822.89'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '11.08'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.660'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.64'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.084.22'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.56'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0792'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '239.22'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.13%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.90'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('E22').Value = '  -1.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.01'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('E24').Value = '  +3.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '172.43'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.83'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.44'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.37%  '
$ws.Range('E28').Value = '  -1.25%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.60'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +20.26%  '
$ws.Range('E30').Value = '  +0.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.332.46'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +37.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.04'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.75%  '
$ws.Range('E33').Value = '  +3.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.99'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.77'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.78%  '
$ws.Range('E36').Value = '  +5.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '92.83'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.682'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.67%  '
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.29'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.308.72'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.998'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '14.63'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('E44').Value = '  +1.37%  '
$ws.Range('E45').Value = '  -5.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.76'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.50%  '
$ws.Range('E47').Value = '  +4.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0512'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.996.64'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.95%  '
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('E51').Value = '  +4.77%  '
